$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# "Arreglo y limpiado de tabla personalizada"
#
# Old layout: a numeric row-code lived in column A (bold/bordered
# "header" style) and the real data started in column B (with the
# "TOTAL" figures in column C, etc.). The fix:
#   * drops that numeric-code column entirely,
#   * shifts every remaining column one place to the left,
#   * drops the old "TOTAL" summary row and re-sources the three
#     data rows (regimen_general / TODOS LOS CENTROS / Hombres),
#   * appends three new descriptive columns: Tabla, Sección, Subsección.
# ------------------------------------------------------------------

# --- 1. Fix up styling first, while the header style still lives on
#        the existing header cells (C1 currently carries it) ---

# Grow the bold/centered/bordered header style onto the new header
# cells (old column A is freed up, three new trailing columns appear).
$ws.Range("C1").Copy()
$ws.Range("A1").PasteSpecial(-4122)
$ws.Range("T1:V1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Old column A (rows 2-4) carried the header-ish style for its numeric
# codes; the new column A is a plain text column, so clear it back to
# the default "Normal" style.
$ws.Range("A2:A4").Style = "Normal"

# --- 2. Header row (row 1) text, now that styles are settled ---
$ws.Cells.Item(1, 1).Value = ""
$ws.Cells.Item(1, 2).Value = "TOTAL"
$ws.Cells.Item(1, 3).Value = "E. Infantil - Primer ciclo"
$ws.Cells.Item(1, 4).Value = "E. Infantil - Segundo ciclo"
$ws.Cells.Item(1, 5).Value = "E. Primaria"
$ws.Cells.Item(1, 6).Value = "Educación Especial"
$ws.Cells.Item(1, 7).Value = "ESO"
$ws.Cells.Item(1, 8).Value = "Bachillerato"
$ws.Cells.Item(1, 9).Value = "Bachillerato a distancia"
$ws.Cells.Item(1, 10).Value = "CF Grado Básico"
$ws.Cells.Item(1, 11).Value = "CF Grado Medio"
$ws.Cells.Item(1, 12).Value = "CF Grado Medio a distancia"
$ws.Cells.Item(1, 13).Value = "Cursos de Especialización Grado Medio"
$ws.Cells.Item(1, 14).Value = "Cursos de Especialización Grado Medio a distancia"
$ws.Cells.Item(1, 15).Value = "CF Grado Superior"
$ws.Cells.Item(1, 16).Value = "CF Grado Superior a distancia"
$ws.Cells.Item(1, 17).Value = "Cursos de Especialización Grado Superior"
$ws.Cells.Item(1, 18).Value = "Cursos de Especialización Grado Superior a distancia"
$ws.Cells.Item(1, 19).Value = "Otros Programas Formativos"
$ws.Cells.Item(1, 20).Value = "Tabla"
$ws.Cells.Item(1, 21).Value = "Sección"
$ws.Cells.Item(1, 22).Value = "Subsección"

# --- 3. Data rows 2-4 ---
# Row 2
$ws.Cells.Item(2, 1).Value = "01 ANDALUCÍA"
$ws.Cells.Item(2, 2).Value = 824964
$ws.Cells.Item(2, 3).Value = 55414
$ws.Cells.Item(2, 4).Value = 113556
$ws.Cells.Item(2, 5).Value = 273456
$ws.Cells.Item(2, 6).Value = 6306
$ws.Cells.Item(2, 7).Value = 212135
$ws.Cells.Item(2, 8).Value = 57897
$ws.Cells.Item(2, 9).Value = 4420
$ws.Cells.Item(2, 10).Value = 10370
$ws.Cells.Item(2, 11).Value = 40959
$ws.Cells.Item(2, 12).Value = 1111
$ws.Cells.Item(2, 13).Value = 130
$ws.Cells.Item(2, 14).Value = 0
$ws.Cells.Item(2, 15).Value = 40668
$ws.Cells.Item(2, 16).Value = 7511
$ws.Cells.Item(2, 17).Value = 712
$ws.Cells.Item(2, 18).Value = 0
$ws.Cells.Item(2, 19).Value = 319
$ws.Cells.Item(2, 20).Value = "regimen_general"
$ws.Cells.Item(2, 21).Value = "TODOS LOS CENTROS"
$ws.Cells.Item(2, 22).Value = "Hombres"

# Row 3
$ws.Cells.Item(3, 1).Value = "01 ANDALUCÍA"
$ws.Cells.Item(3, 2).Value = 778862
$ws.Cells.Item(3, 3).Value = 51205
$ws.Cells.Item(3, 4).Value = 107767
$ws.Cells.Item(3, 5).Value = 257765
$ws.Cells.Item(3, 6).Value = 3120
$ws.Cells.Item(3, 7).Value = 199797
$ws.Cells.Item(3, 8).Value = 67368
$ws.Cells.Item(3, 9).Value = 4117
$ws.Cells.Item(3, 10).Value = 4112
$ws.Cells.Item(3, 11).Value = 32521
$ws.Cells.Item(3, 12).Value = 2231
$ws.Cells.Item(3, 13).Value = 24
$ws.Cells.Item(3, 14).Value = 0
$ws.Cells.Item(3, 15).Value = 38411
$ws.Cells.Item(3, 16).Value = 10064
$ws.Cells.Item(3, 17).Value = 176
$ws.Cells.Item(3, 18).Value = 0
$ws.Cells.Item(3, 19).Value = 184
$ws.Cells.Item(3, 20).Value = "regimen_general"
$ws.Cells.Item(3, 21).Value = "TODOS LOS CENTROS"
$ws.Cells.Item(3, 22).Value = "Hombres"

# Row 4
$ws.Cells.Item(4, 1).Value = "01 ANDALUCÍA"
$ws.Cells.Item(4, 2).Value = 1155976
$ws.Cells.Item(4, 3).Value = 38467
$ws.Cells.Item(4, 4).Value = 166760
$ws.Cells.Item(4, 5).Value = 401713
$ws.Cells.Item(4, 6).Value = 6770
$ws.Cells.Item(4, 7).Value = 312903
$ws.Cells.Item(4, 8).Value = 99284
$ws.Cells.Item(4, 9).Value = 8537
$ws.Cells.Item(4, 10).Value = 11010
$ws.Cells.Item(4, 11).Value = 48284
$ws.Cells.Item(4, 12).Value = 1260
$ws.Cells.Item(4, 13).Value = 154
$ws.Cells.Item(4, 14).Value = 0
$ws.Cells.Item(4, 15).Value = 51219
$ws.Cells.Item(4, 16).Value = 8238
$ws.Cells.Item(4, 17).Value = 874
$ws.Cells.Item(4, 18).Value = 0
$ws.Cells.Item(4, 19).Value = 503
$ws.Cells.Item(4, 20).Value = "regimen_general"
$ws.Cells.Item(4, 21).Value = "TODOS LOS CENTROS"
$ws.Cells.Item(4, 22).Value = "Hombres"

